$d = $word.ActiveDocument

$range = $d.Content
$range.Find.Execute(
    "functionality of all links",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "functionality of AMPscript code, links",
    2)
